$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=f8d0ee9e9b967d7b0cf389c9a3d9dc7293000759d60cac4cb40a579a5daf5c1a'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=138f1fc4cfc2ae80f53053f2b409a2e27ca0dbaf33124bee1f0fa66c373c78ef'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=1dc6ed4d262f4aa6221c7748023501fe7adc73a96846bf240f1c98c96794e1ed'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=d7ca24a03c1c1c75cfc37454ebe7c21dbf981da9f505c18f625faa341eed203f'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T000004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2008%3A00%3A04%20GMT&X-Amz-Signature=c755dd9c3672acf979ca38e550f5ccca9027a161506ef934949f605bc9426bbd'

$wb.Save()